# Add notify for all admin view
# Updates the "malophp" (class code) column to encode the related
# subject code (mamh) together with a per-subject class sequence number,
# e.g. "LHP001" -> "MH001-LHP001".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newCodes = @{
    2  = "MH001-LHP001"
    3  = "MH001-LHP002"
    4  = "MH003-LHP001"
    5  = "MH003-LHP002"
    6  = "MH005-LHP001"
    7  = "MH005-LHP002"
    8  = "MH007-LHP001"
    9  = "MH007-LHP002"
    10 = "MH009-LHP001"
    11 = "MH009-LHP002"
    12 = "MH011-LHP001"
    13 = "MH011-LHP002"
    14 = "MH012-LHP001"
}

foreach ($row in $newCodes.Keys) {
    $ws.Cells.Item($row, 1).Value = $newCodes[$row]
}

$ws.Range("A15").Select()
